# 11/7 updates.. a lot
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# D23: new value "codeup : 6017-6020", formatted like C23 (fill style)
# ---------------------------------------------------------------------
$ws.Range("D23").Value = "codeup : 6017-6020"

# ---------------------------------------------------------------------
# Row 27: shift old C27 ("String Split and Join") into D27 (keeping its
# original highlighted format), then turn B27/C27 into "X" cells that
# match the plain style used by A27.
# ---------------------------------------------------------------------
$oldC27Value = $ws.Range("C27").Value2

# Preserve C27's current fill/format on D27 before C27's format changes.
$ws.Range("C27").Copy()
$ws.Range("D27").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

$ws.Range("B27").Value = "X"
$ws.Range("C27").Value = "X"
$ws.Range("D27").Value = $oldC27Value

# ---------------------------------------------------------------------
# D28: new value "codeup : 4012 (석차계산)", formatted like C28
# ---------------------------------------------------------------------
$ws.Range("D28").Value = "codeup : 4012 (석차계산)"

# ---------------------------------------------------------------------
# D4: new value "4012(석차계산)" (cell already formatted, style preserved)
# ---------------------------------------------------------------------
$ws.Range("D4").Value = "4012(석차계산)"

# ---------------------------------------------------------------------
# D25: new value "codeup : 1025-1030", formatted like C25
# ---------------------------------------------------------------------
$ws.Range("D25").Value = "codeup : 1025-1030"

# ---------------------------------------------------------------------
# Formatting fix-ups (copy format only, values already set above)
# ---------------------------------------------------------------------
$ws.Range("C23").Copy()
$ws.Range("D23").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

$ws.Range("C25").Copy()
$ws.Range("D25").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

$ws.Range("C28").Copy()
$ws.Range("D28").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

$ws.Range("A27").Copy()
$ws.Range("B27").PasteSpecial($xlPasteFormats)
$ws.Range("C27").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# sheetView: scroll/top-left cell moved to A10, selection moved to F23
# ---------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F23").Select()
